$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper: copy the number/cell format from a source cell to a destination
# cell (format only - xlPasteFormats), leaving the destination's value
# untouched (it will be set separately afterwards, if needed).
# ---------------------------------------------------------------------------
function Copy-Format($srcRow, $srcCol, $dstRow, $dstCol) {
    $ws.Cells.Item($srcRow, $srcCol).Copy() | Out-Null
    $ws.Cells.Item($dstRow, $dstCol).PasteSpecial(-4122) | Out-Null
}

# ---------------------------------------------------------------------------
# 1. Text cells that pick up word-wrap (style moves from the "no-wrap" xf
#    to the already-existing "wrap" xf - same font, just wrapText=1).
# ---------------------------------------------------------------------------
$ws.Cells.Item(4, 3).WrapText = $true    # C4
$ws.Cells.Item(5, 2).WrapText = $true    # B5
$ws.Cells.Item(5, 3).WrapText = $true    # C5
$ws.Cells.Item(6, 2).WrapText = $true    # B6
$ws.Cells.Item(6, 3).WrapText = $true    # C6
$ws.Cells.Item(9, 2).WrapText = $true    # B9

# ---------------------------------------------------------------------------
# 2. Row 3 - extend the year header row with 2021 / 2022 / 2023.
# ---------------------------------------------------------------------------
Copy-Format 3 4 3 15            # O3 <- D3 format
$ws.Cells.Item(3, 15).Value = 2021
Copy-Format 3 4 3 16            # P3 <- D3 format
$ws.Cells.Item(3, 16).Value = 2022
Copy-Format 3 4 3 17            # Q3 <- D3 format
$ws.Cells.Item(3, 17).Value = 2023

# ---------------------------------------------------------------------------
# 3. Data rows 4-11 : fill column N (2020) where it was still blank, and add
#    the new columns O (2021), P (2022) and Q (2023).
#    Column D of the same row always already carries the exact target
#    number style, so it is used as the format donor.
# ---------------------------------------------------------------------------

# --- Row 4 ------------------------------------------------------------
Copy-Format 4 4 4 14 ; $ws.Cells.Item(4, 14).Value = 92.9
Copy-Format 4 4 4 15 ; $ws.Cells.Item(4, 15).Value = 105.5
Copy-Format 4 4 4 16 ; $ws.Cells.Item(4, 16).Value = 109
Copy-Format 4 4 4 17 ; $ws.Cells.Item(4, 17).Value = 106.2

# --- Row 5 ------------------------------------------------------------
Copy-Format 5 4 5 14 ; $ws.Cells.Item(5, 14).Value = 106.3
Copy-Format 5 4 5 15 ; $ws.Cells.Item(5, 15).Value = 111.90503981851454
Copy-Format 5 4 5 16 ; $ws.Cells.Item(5, 16).Value = 113.92290931741762
Copy-Format 5 4 5 17 ; $ws.Cells.Item(5, 17).Value = 110.8

# --- Row 6 (also M6 value correction 104.1 -> 104.26) ------------------
$ws.Cells.Item(6, 13).Value = 104.26
Copy-Format 6 4 6 14 ; $ws.Cells.Item(6, 14).Value = 121.27
Copy-Format 6 4 6 15 ; $ws.Cells.Item(6, 15).Value = 111.5
Copy-Format 6 4 6 16 ; $ws.Cells.Item(6, 16).Value = 105.1
Copy-Format 6 4 6 17 ; $ws.Cells.Item(6, 17).Value = 109.3

# --- Row 7 (Q7 stays blank, formatted only) -----------------------------
Copy-Format 7 4 7 14 ; $ws.Cells.Item(7, 14).Value = -19734.0366
Copy-Format 7 4 7 15 ; $ws.Cells.Item(7, 15).Value = -1763.6
Copy-Format 7 4 7 16 ; $ws.Cells.Item(7, 16).Value = -10400.700000000001
Copy-Format 7 4 7 17                                   # Q7 blank

# --- Row 8 (M8 value filled in, Q8 stays blank) -------------------------
$ws.Cells.Item(8, 13).Value = 319474.59999999998
Copy-Format 8 4 8 14 ; $ws.Cells.Item(8, 14).Value = 407116.85000000003
Copy-Format 8 4 8 15 ; $ws.Cells.Item(8, 15).Value = 436586.8
Copy-Format 8 4 8 16 ; $ws.Cells.Item(8, 16).Value = 477967.8
Copy-Format 8 4 8 17                                   # Q8 blank

# --- Row 9 (row grows taller, Q9 stays blank) ---------------------------
Copy-Format 9 4 9 14 ; $ws.Cells.Item(9, 14).Value = 81.599999999999994
Copy-Format 9 4 9 15 ; $ws.Cells.Item(9, 15).Value = 146.4
Copy-Format 9 4 9 16 ; $ws.Cells.Item(9, 16).Value = 144.69999999999999
Copy-Format 9 4 9 17                                   # Q9 blank
$ws.Rows.Item(9).RowHeight = 24

# --- Row 10 (Q10 stays blank) -------------------------------------------
Copy-Format 10 4 10 14 ; $ws.Cells.Item(10, 14).Value = 1973.2
Copy-Format 10 4 10 15 ; $ws.Cells.Item(10, 15).Value = 2752.1
Copy-Format 10 4 10 16 ; $ws.Cells.Item(10, 16).Value = 2254.6999999999998
Copy-Format 10 4 10 17                                   # Q10 blank

# --- Row 11 (Q11 stays blank) -------------------------------------------
Copy-Format 11 4 11 14 ; $ws.Cells.Item(11, 14).Value = 3718.8
Copy-Format 11 4 11 15 ; $ws.Cells.Item(11, 15).Value = 5580.2
Copy-Format 11 4 11 16 ; $ws.Cells.Item(11, 16).Value = 9803.2000000000007
Copy-Format 11 4 11 17                                   # Q11 blank

# ---------------------------------------------------------------------------
# 4. Columns A:C become a single uniform width.
# ---------------------------------------------------------------------------
$ws.Range("A:C").ColumnWidth = 36.65

# ---------------------------------------------------------------------------
# 5. Reset the selection back to the top-left cell (drops the stored
#    I19 selection that used to be saved with the sheet).
# ---------------------------------------------------------------------------
$ws.Range("A1").Select()

Write-Host "edit complete"
